# Generate Report for Handoff
# Update the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# for the last row (5647227f-a389-4afe-a183-5cd2c0bab255) on each sheet, reflecting
# a newly generated handoff report.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date", row 7
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-09-03 12:45:35"

# zh-cn sheet: column H = "Latest Handoff Datetime", row 7
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-09-03 12:45:30"

# de-de sheet: column H = "Latest Handoff Datetime", row 7
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-09-03 12:45:35"
